$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 29 (old rows 29-35 shift down to 30-36).
$ws.Rows("29:29").Insert()

# 2. Renumber the "Number" column for the rows that were pushed down by the insert.
$ws.Range("A29").Value = 28
$ws.Range("A30").Value = 29
$ws.Range("A31").Value = 30
$ws.Range("A32").Value = 31
$ws.Range("A33").Value = 32
$ws.Range("A34").Value = 33
$ws.Range("A35").Value = 34
$ws.Range("A36").Value = 35
$ws.Range("A37").Value = 36
$ws.Range("A38").Value = 37

# 3. Variable names (column B) for the new/appended rows.
$ws.Range("B29").Value = "VFA_"
$ws.Range("B37").Value = "Underweight"
$ws.Range("B38").Value = "Agegroup"

# 4. Variable labels (column C).
$ws.Range("C29").Value = "Visceral fat area (cm^2)"
$ws.Range("C37").Value = "Underweight (BMI <18.5 kg/m^2)"
$ws.Range("C38").Value = "Age group (years)"

# 5. Variable type (column E).
$ws.Range("E29").Value = "Numeric"
$ws.Range("E37").Value = "Character"
$ws.Range("E38").Value = "Character"

# 6. Coded responses (column D), wrapped text.
$ws.Range("D29").Value = "0=No`n1=Yes"
$ws.Range("D31").Value = "0=No`n1=Yes"
$ws.Range("D37").Value = "0=No`n1=Yes"
$ws.Range("D38").Value = "1=20-29 years`n2=30-39 years`n3=40-49 years`n4=50-59 years`n5=60-69 years`n6=70 years and older"

$ws.Range("D29").WrapText = $true
$ws.Range("D31").WrapText = $true
$ws.Range("D37").WrapText = $true
$ws.Range("D38").WrapText = $true

$ws.Rows("29:29").RowHeight = 29
$ws.Rows("31:31").RowHeight = 29
$ws.Rows("37:37").RowHeight = 29
$ws.Rows("38:38").RowHeight = 87

# 7. Update the sheet view selection to match the edited state (the engine
#    does not round-trip topLeftCell scroll position).
$ws.Range("E39").Select()
